$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: scroll/selection moved from I10 area to N10 area ---
$ws.Application.ActiveWindow.TopLeftCell = $ws.Range("G6")
$ws.Range("N10:N13").Select()

# --- Row 6: add hyperlink to evidence screenshot in M6 (Evidencia column) ---
$ws.Hyperlinks.Add($ws.Range("M6"), "https://github.com/Goriguen/StockFreezer/blob/10f2d977b450b81e77b3299950c65c4f9036fc87/docs/qa_testing/evidencia_screenshots/bug_cp005_menu_codigo.jpg")

# --- Row 10: QA case CP-005 now passes ---
# Estado: FALLÓ -> PASÓ
$ws.Range("J10").Value = "PASÓ"

# Fecha Corrección: now carries the same date as Fecha Prueba (K10), using
# the same date style already used by the Fecha Prueba / Fecha Corrección
# columns elsewhere in the sheet.
$ws.Range("K10").Copy()
$ws.Range("L10:L13").PasteSpecial(-4122)
$ws.Range("L10").Value2 = $ws.Range("K10").Value2

# Evidencia (M10) adopts the plain text style used by the Comentarios column
$ws.Range("N10").Copy()
$ws.Range("M10:M13").PasteSpecial(-4122)
$ws.Range("M10").Value = "/"

# Comentarios: bug got fixed
$ws.Range("N10").Value = "El sistema eliminó perfectamente el producto en las coordenadas indicadas."

$ws.Application.CutCopyMode = $false
